# Updated ToT-ToA, Assessor-Trainer ans TP-Candididate Workflow test cases.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Candidate full-name values renamed (Bulk -> BulkTR naming) ---
$ws.Range("B5").Value = "PalsVpBulkTRA"
$ws.Range("B6").Value = "PalsVpBulkTRB"
$ws.Range("B7").Value = "PalsVpBulkTRC"
$ws.Range("B8").Value = "PalsVpBulkTRD"

# --- Place of birth swaps ---
$ws.Range("E5").Value = "Aimangala"
$ws.Range("E7").Value = "Kasaragod"
$ws.Range("E8").Value = "Big Lapati"

# --- Row 5: Current Address block now points to Tirunelveli, Tamil Nadu ---
$ws.Range("T5").Value = "#401-406, World Mark 1, West Wing, Aerocity"
$ws.Range("U5").Value = "Tirunelveli"
$ws.Range("V5").Value = "TIRUNELVELI"
$ws.Range("W5").Value = "'677558"
$ws.Range("X5").Value = "TAMIL NADU"
$ws.Range("Y5").Value = "Tirunelveli"

# New pasted-style font (Arial 10, #222222) applied to the cells that came
# from the pasted web content.
foreach ($addr in @("U5", "V5", "X5", "Y5")) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.Color = 2236962
}

# --- View state: scroll back to the top-left and select B6 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select()
